$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '27.551.76'
$ws.Cells.Item(2, 5).Value = '  -0.13%  '
$ws.Cells.Item(3, 4).Value = '1.756.24'
$ws.Cells.Item(3, 5).Value = '  +0.13%  '
$ws.Cells.Item(4, 5).Value = '  +0.04%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '324.44'
$ws.Cells.Item(5, 5).Value = '  +0.00%  '
$ws.Cells.Item(6, 5).Value = '  +0.04%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.4564'
$ws.Cells.Item(7, 5).Value = '  +1.55%  '
$ws.Cells.Item(8, 5).Value = '  -1.86%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.07473'
$ws.Cells.Item(9, 5).Value = '  -0.52%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '41.54'
$ws.Cells.Item(10, 5).Value = '  -1.59%  '
$ws.Cells.Item(11, 5).Value = '  -1.71%  '
$ws.Cells.Item(12, 5).Value = '  +0.04%  '
$ws.Cells.Item(13, 5).Value = '  +0.41%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '6.007'
$ws.Cells.Item(14, 5).Value = '  -0.73%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '7.155'
$ws.Cells.Item(15, 5).Value = '  -0.40%  '
$ws.Cells.Item(16, 4).Value = '1.752.88'
$ws.Cells.Item(16, 5).Value = '  -0.01%  '
$ws.Cells.Item(17, 5).Value = '  +0.83%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.00001054'
$ws.Cells.Item(18, 5).Value = '  -1.17%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.06391'
$ws.Cells.Item(19, 5).Value = '  -0.04%  '
$ws.Cells.Item(20, 5).Value = '  -0.04%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '17.08'
$ws.Cells.Item(21, 5).Value = '  +0.77%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '5.740'
$ws.Cells.Item(22, 5).Value = '  -2.21%  '
$ws.Cells.Item(23, 4).Value = '27.596.48'
$ws.Cells.Item(23, 5).Value = '  -0.06%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '11.23'
$ws.Cells.Item(24, 5).Value = '  -0.04%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.071'
$ws.Cells.Item(25, 5).Value = '  -1.88%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '165.59'
$ws.Cells.Item(26, 5).Value = '  +2.53%  '
$ws.Cells.Item(27, 5).Value = '  -1.45%  '
$ws.Cells.Item(28, 4).Value = '1.955.27'
$ws.Cells.Item(28, 5).Value = '  +0.06%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '2.140'
$ws.Cells.Item(29, 5).Value = '  +0.77%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '125.59'
$ws.Cells.Item(30, 5).Value = '  +0.14%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.084'
$ws.Cells.Item(31, 5).Value = '  -0.47%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.09236'
$ws.Cells.Item(32, 5).Value = '  +2.14%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '3.657'
$ws.Cells.Item(33, 5).Value = '  +0.58%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '5.527'
$ws.Cells.Item(34, 5).Value = '  -0.77%  '
$ws.Cells.Item(35, 5).Value = '  -2.30%  '
$ws.Cells.Item(36, 5).Value = '  -1.24%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.2093'
$ws.Cells.Item(37, 5).Value = '  +0.24%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.06022'
$ws.Cells.Item(38, 5).Value = '  +0.65%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.6284'
$ws.Cells.Item(39, 5).Value = '  -1.50%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '4.926'
$ws.Cells.Item(40, 5).Value = '  -1.08%  '
$ws.Cells.Item(41, 5).Value = '  -1.53%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.387'
$ws.Cells.Item(42, 5).Value = '  +0.06%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '7.807'
$ws.Cells.Item(43, 5).Value = '  +0.02%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '13.23'
$ws.Cells.Item(44, 5).Value = '  -0.40%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '3.718'
$ws.Cells.Item(45, 5).Value = '  +0.14%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.5869'
$ws.Cells.Item(46, 5).Value = '  -0.45%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '122.35'
$ws.Cells.Item(47, 5).Value = '  +0.61%  '
$ws.Cells.Item(48, 5).Value = '  -1.12%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.06896'
$ws.Cells.Item(49, 5).Value = '  +0.38%  '
$ws.Cells.Item(50, 5).Value = '  -2.82%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '72.28'
$ws.Cells.Item(51, 5).Value = '  -0.33%  '
